$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Shift B and C columns up by 4 rows: new row r (2..16) gets old row (r+4) values
    for ($r = 2; $r -le 16; $r++) {
        $srcRow = $r + 4
        $bVal = $ws.Cells.Item($srcRow, 2).Value()
        $cVal = $ws.Cells.Item($srcRow, 3).Value()
        $ws.Cells.Item($r, 2).Value = $bVal
        $ws.Cells.Item($r, 3).Value = $cVal
    }

    # Delete now-obsolete rows 17-20
    $ws.Range("A17:C20").EntireRow.Delete()
}
